$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '54.363.34'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.284.42'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '500.76'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.10%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '129.42'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.10%  '
$ws.Range('E6').ClearFormats()
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E7').ClearFormats()
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0954'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.45%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.152'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('E10').ClearFormats()
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +3.20%  '
$ws.Range('E11').ClearFormats()
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.689.13'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('E13').ClearFormats()
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.96'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +6.51%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '54.299.51'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('E15').ClearFormats()
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.303.35'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.26'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +3.23%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.13'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.13%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '304.86'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.80%  '
$ws.Range('E20').ClearFormats()
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +1.96%  '
$ws.Range('E21').ClearFormats()
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '62.02'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -2.41%  '
$ws.Range('E23').ClearFormats()
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('E24').ClearFormats()
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.34'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '173.84'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +6.66%  '
$ws.Range('E27').ClearFormats()
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('E28').ClearFormats()
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.99'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.79%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0691'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.09'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.43%  '
$ws.Range('E31').ClearFormats()
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('E32').ClearFormats()
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '17.82'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +1.74%  '
$ws.Range('E33').ClearFormats()
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E34').ClearFormats()
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.937'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +8.64%  '
$ws.Range('E35').ClearFormats()
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.44%  '
$ws.Range('E36').ClearFormats()
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.77'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('E37').ClearFormats()
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.60%  '
$ws.Range('E38').ClearFormats()
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.40'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('E40').ClearFormats()
$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('B41').ClearFormats()
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C41').ClearFormats()
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.01'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +1.68%  '
$ws.Range('E41').ClearFormats()
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'Aave'
$ws.Range('B42').ClearFormats()
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C42').ClearFormats()
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '124.92'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('E42').ClearFormats()
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0496'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +3.20%  '
$ws.Range('E43').ClearFormats()
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('E44').ClearFormats()
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.549'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '240.90'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('E46').ClearFormats()
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.67%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0207'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('E48').ClearFormats()
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.07%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '16.40'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('E50').ClearFormats()
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.46%  '
$ws.Range('E51').ClearFormats()
